$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Reference style (currency number format) used by the existing price column
$priceFormat = $ws.Range("C20").NumberFormat

$data = @(
    @(21, "Aloo Patty", 20, "Aloo Patty.jpg"),
    @(22, "Paneer Patty", 25, "Paneer Patty.jpg"),
    @(23, "Butter Patty", 30, "Butter Patty.jpg"),
    @(24, "Pastry Pineapple", 25, "Pastry Pineappl.jpg"),
    @(25, "Pastry Chocolate", 45, "Pastry Chocolate.jpg"),
    @(26, "Fanta20", 10, "Fanta20.jpg"),
    @(27, "Thumsup20", 20, "Thumsup20.jpg"),
    @(28, "Frooti20", 20, "Frooti20jpg"),
    @(29, "Maza10", 10, "Maza-Tetra-pack.jpg"),
    @(30, "Salted French Frise", 80, "Salted French Fires image.jpg"),
    @(31, "Peri-Peri French Frise", 100, "Peri- Peri French Fires.jpg")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("C$r").NumberFormat = $priceFormat
    $ws.Range("D$r").Value = $row[3]
}

$ws.Range("D26").Select()
